# Adds a new "as of" forecast column (AI) and a new observation row (47)
# for 2020-05-20 to both the "cases" and "deaths" sheets, and backfills
# the newly-observed value for 2020-05-06 (row 33, column B).
#
# Column AI header date:  2020-05-06  (reuses existing shared string)
# New row 47 date:        2020-05-20  (brand new shared string)

$wb = $excel.ActiveWorkbook

# ---- helpers -------------------------------------------------------

# Write a literal text value into a cell without Excel's "looks like a
# date" auto-conversion kicking in (dates in this table are stored as
# plain shared-string text, e.g. "2020-05-06").
function Set-TextCell {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Write a numeric value into a cell.
function Set-NumberCell {
    param($range, $number)
    $range.Value2 = $number
}

# Force a cell to be persisted as an explicit-but-empty cell (matches
# the sparse "<c .../>" placeholders already used throughout the sheet).
function Set-EmptyCell {
    param($range)
    $range.Style = "Normal"
}

# ---- per-sheet edits -------------------------------------------------

function Update-ForecastSheet {
    param($ws, $b33Value, $aiValues)

    # New header cell AI1: forecast "as of" date for the new column.
    Set-TextCell $ws.Range("AI1") "2020-05-06"

    # Give every existing data row (2-33) an empty AI cell, matching the
    # sparse layout used by every other column.
    $ws.Range("AI2:AI33").Style = "Normal"

    # Backfill the now-observed value for row 33 (date 2020-05-06).
    Set-NumberCell $ws.Range("B33") $b33Value

    # Fill in the forecast values for rows 34-46 in the new AI column.
    foreach ($rowNum in 34..46) {
        Set-NumberCell $ws.Cells.Item($rowNum, 35) $aiValues[$rowNum]
    }

    # New row 47 (date 2020-05-20): create the empty cells A47:AH47 ...
    $ws.Range("A47:AH47").Style = "Normal"
    # ... set its date label ...
    Set-TextCell $ws.Range("A47") "2020-05-20"
    # ... and its single forecast value in the new AI column.
    Set-NumberCell $ws.Range("AI47") $aiValues[47]
}

$wsCases = $wb.Worksheets.Item("cases")
$casesAI = @{
    34 = 134005
    35 = 142591
    36 = 150754
    37 = 157802
    38 = 164903
    39 = 172253
    40 = 178227
    41 = 184377
    42 = 189342
    43 = 194518
    44 = 199012
    45 = 203603
    46 = 208168
    47 = 212154
}
Update-ForecastSheet $wsCases 125218 $casesAI

$wsDeaths = $wb.Worksheets.Item("deaths")
$deathsAI = @{
    34 = 9033
    35 = 9520
    36 = 9971
    37 = 10352
    38 = 10738
    39 = 11134
    40 = 11446
    41 = 11773
    42 = 12032
    43 = 12308
    44 = 12543
    45 = 12789
    46 = 13037
    47 = 13253
}
Update-ForecastSheet $wsDeaths 8536 $deathsAI
